# MIDI Mappings Freeverb.xlsx - revise MIDI mappings / documentation for
# the stand-alone effects table on the "By Module" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constants
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138

# --- Update the mapped values -------------------------------------------
# Row 3 becomes the "Enable" switch mapped to SW-4 (CC 104)
$ws.Range("A3").Value = "Enable"
$ws.Range("B3").Value = 104
$ws.Range("D3").Value = "switch"

# Row 4 becomes "Wet Dry" mapped to HADC-0 (CC 2)
$ws.Range("A4").Value = "Wet Dry"
$ws.Range("B4").Value = 2
$ws.Range("D4").Value = "knob"

# Row 5 becomes "Damp" mapped to HADC-1 (CC 3)
$ws.Range("A5").Value = "Damp"
$ws.Range("B5").Value = 3

# Row 6 becomes "Room Size" mapped to HADC-2 (CC 4)
$ws.Range("A6").Value = "Room Size"
$ws.Range("B6").Value = 4

# --- Thicken the outer border of the table -------------------------------
# Top edge of row 4 (A4, B4, D4) becomes a thin line (previously unset) -
# apply this before the right-edge pass below so that D4 only needs a
# single additional (new) border combination instead of two.
foreach ($addr in @("A4","B4","D4")) {
    $b = $ws.Range($addr).Borders.Item($xlEdgeTop)
    $b.LineStyle = $xlContinuous
    $b.Weight = $xlThin
}

# Right edge of column D (the rightmost column of the table) becomes medium
foreach ($addr in @("D1","D2","D3","D4","D5","D6")) {
    $b = $ws.Range($addr).Borders.Item($xlEdgeRight)
    $b.LineStyle = $xlContinuous
    $b.Weight = $xlMedium
}

# Left edge of row 2 (column A) becomes medium
$bA2 = $ws.Range("A2").Borders.Item($xlEdgeLeft)
$bA2.LineStyle = $xlContinuous
$bA2.Weight = $xlMedium

# --- Update the active selection -----------------------------------------
$ws.Range("A20").Select() | Out-Null
